# Created clock panel config and scroll wheel for Ti setup panel.
#
# 1. Typography sheet: the "Default" typography's Wildcard Ranges
#    (column I, rows 4-8) gain lowercase/uppercase letters alongside digits
#    so the new alphanumeric text fields can render.
# 2. Translation sheet: a new translation row (32) is added for the Ti
#    setup panel's scroll-wheel placeholder text ("<value>").

$wb  = $excel.ActiveWorkbook
$typography  = $wb.Worksheets.Item("Typography")
$translation = $wb.Worksheets.Item("Translation")

# --- Typography sheet: widen the wildcard range for the Default font ---
$typography.Range("I4:I8").Value = "a-z,A-Z,0-9"

# --- Translation sheet: new row for the scroll wheel's current value ---
$translation.Range("B32").Value = "SingleUseId29"
$translation.Range("C32").Value = "Default"
$translation.Range("D32").Value = "Center"
$translation.Range("E32").Value = "LTR"
$translation.Range("F32").Value = "<value>"
